# Apply edits to "Other Info" worksheet (sheet2): add Program Content section
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Row 1 / Row 2 / Row 3 text values, set in the exact order the original
# authors entered them so new shared-string indices line up with the target. ---
$ws.Range("H1").Value = "Program Content :                                                      "
$ws.Range("H2").Value = "Science"
$ws.Range("I2").Value = "Socail Science"
$ws.Range("J2").Value = "Food & Agriculture"
$ws.Range("K2").Value = "Forestry"
$ws.Range("L2").Value = "Sustainability/Enviromental Ed"
$ws.Range("M2").Value = "Education Arts"
$ws.Range("N2").Value = "Language Arts"
$ws.Range("O2").Value = "Math"
$ws.Range("P2").Value = "Geography"
$ws.Range("Q2").Value = "STEM/STEAM"
$ws.Range("R2").Value = "Visual and Performing Arts"
$ws.Range("S2").Value = "Physical/Health Ed"
$ws.Range("T2").Value = "Other"
$ws.Range("U1").Value = "Please describe how your ODS program’s content (questions 3 and 4) is integrated with local school curricula in a manner that assists students in meeting state standards related to science, technology, engineering and mathematics (STEM), and the Next Generation Science Standards."
$ws.Range("U3").Value = "MOSS Curriculum is rooted in Next Generation Science Standards. Each week, teachers are able to choose a Disciplinary Core Idea (relevent to the local ecosystem) on which to focus the lessons throughout the week. Insturctors deliver the content through the Practices of Engineering and Science during the program, and guide students through inquiry-based science investigations in Pondeosa State Park. Each week, teachers may also choose between an `"Inquiry Project`" or a `"Community Engineering Project`" for the students to focus on for an entire day at the end of the week. Students design their own projects, in which they either ask a question they attempt to answer by collecting data in the natural enviornment, or identify a natural resources problem in the local community that they can attempt to solve using the content they learn throughout the week. "
$ws.Range("V1").Value = "Please describe if (and if so, how) your ODS program is offered in a bilingual format."
$ws.Range("V3").Value = "MOSS relies on schools to bring translators or interpreters with their multi-lingual students. "
$ws.Range("W1").Value = "Please describe how the program addresses the inequity of outdoor educational opportunities for underserved children in this state."
$ws.Range("W3").Value = "MOSS actively seeks grant funding for students who are unable to attend programs, and supports classrooms that are challenged to find funding through assistance in grant writing and in-kind donations of instrcution by staff and volunteers. "
$ws.Range("X1").Value = "Please describe how your ODS program provides students with opportunities to learn about the interdependence of urban and rural areas."
$ws.Range("X3").Value = "MOSS is located in rural central Idaho. Due to the location, we are constantly seeking connection with student from urban and rural areas to identify unique natural resources issues in their area, and be able to identify the ecological principles that may apply to their hometown. By understanding the basics of ecoystem interactions, students will be better suited to identify issues in their home area, and see that the human landscape is intricately interwoven with the natural landscape. "
$ws.Range("Y1").Value = "What instructional strategies are used during the program (select all that apply.)"
$ws.Range("Y2").Value = "Project based Learning"
$ws.Range("Y3").Value = "X"
$ws.Range("Z2").Value = "Cooperative learning stategies"
$ws.Range("AA2").Value = "Service Learning"
$ws.Range("AB2").Value = "Interdisciplinary instruction"
$ws.Range("AC2").Value = "Inquiry-based instruction"
$ws.Range("AD2").Value = "Social Emotional learning"
$ws.Range("AE2").Value = "Socio scientific issues"
$ws.Range("AF2").Value = "Other (list)"
$ws.Range("AF3").Value = "Place-Based Learning"
$ws.Range("AG1").Value = "Please describe how learning is extended back to the classroom and annual curriculum. (This may include but is not limited to pre-and post-activities, assessments of learning, teacher professional development, etc.)"
$ws.Range("AG3").Value = "MOSS provides pre and post-visit activities for teachers to utilize in the classroom upon request. The K12 Programs Coordinator diligently seeks conversations with teachers before their arrival at MOSS to understand how their experience will apply to what the students are learning in the classroom, and how the field experiences can capitalize on projects in a variety of classes (science, math, language arts, social studies, art, etc). "
$ws.Range("AH1").Value = "Please describe how you involve partners in the program. (Partners can include but are not limited to community members, volunteers, government agencies, local non-profits, etc.)"
$ws.Range("AH3").Value = "MOSS programs rely heavily on the local McCall community for assistance in designing place-based curriculum that spans socio-ecological systems in the area. We partner with Idaho State Parks for a location in which to teach; The US Forest Service often provides guest speakers for Evening Programs to discuss careers in natural resources (wildlind fire, forestry, hydrology, soil science, and others); The Nez Perce Tribe provides guest speakers about fisheries and native culture; local recreation industry representatives have provided guest speakers to discuss the effects of use on the natural world and their industry; local ski resorts and the Payette Avalanche Center provide insturction on snow science and backcountry rescue. "
$ws.Range("AI1").Value = "Please briefly describe your success for each of the items below. Use `"n/a`" for any items not addressed in your ODS program."
$ws.Range("AI2").Value = "Higher scores on standardized measures of academic achievement in reading, writing, math, science and social studies."
$ws.Range("AI3").Value = "We have seen a greater understanding and an overall higher performance on standardized tests."
$ws.Range("AJ2").Value = "Greater self-sufficiency and leadership skills"
$ws.Range("AJ3").Value = "Teachers report that students are more self-reliant throughout and after programming due to the residential nature of the program. "
$ws.Range("AK2").Value = "Fewer discipline and classroom management problems"
$ws.Range("AK3").Value = "Students who are identified as having possible behavior issues are rarely seen to exhibit these behaviors in the field. These students are able to channel energy into the field studies"
$ws.Range("AL2").Value = "Increased student engagement and pride in accomplishments"
$ws.Range("AL3").Value = "Students develop, research, and present their own Inquiry and Engineering Projects each week, providing them with pride in starting and completing a project largely on their own as a team. "
$ws.Range("AM2").Value = "Greater proficiency in solving problems and thinking strategically."
$ws.Range("AM3").Value = "Students develop, research, and present their own Inquiry and Engineering Projects each week, providing them with pride in starting and completing a project largely on their own as a team. "
$ws.Range("AN2").Value = "Better application of systems thinking and increased ability to think creatively."
$ws.Range("AN3").Value = "Students are also expected to make connections between the ecosystems they are experiencing and their home environments. "
$ws.Range("AO2").Value = "Improved communication skills and enhanced ability to work in group settings."
$ws.Range("AO3").Value = "Students work collaboratively all thoughout the week. One of the `"goals`" of MOSS is to work collabortively and solve problems as a team, through which they develop a sense of community. Students are also living and working together all week long, which requires them to work together to acheive group goals and adhere to group values identified at the beginning of each week."
$ws.Range("AP2").Value = "Greater enthusiasm for language arts, math, science and social studies"
$ws.Range("AP3").Value = "By the end of each week, students identify themselves as scientists! (Not just a man in a laboratory wearing a labcoat)"
$ws.Range("AQ2").Value = "Increased knowledge and understanding of science content, concepts and processes."
$ws.Range("AQ3").Value = "Students identify themselves as scientists by the end of the week. Throughout the experience, students practice science through a lens of inquiry and discovering connections in the ecosystem. Students are able to think scientifically and identify key aspects of the ecosystem they study. "
$ws.Range("AR2").Value = "Better ability to apply science and civic processes to real-world situations"
$ws.Range("AR3").Value = "Students are asked to include a `"why should we care?`" section of each Inquiry or Engineering project they complete. In this section, students identify that the principles and concepts they implemented in their self-driven projects have greater implications if they are scaled outward in the community or region. "
$ws.Range("AS2").Value = "Improved understanding of mathematical concepts and mastery of math skills."
$ws.Range("AS3").Value = "Students are asked to graph data they collect each week and interpret the graphs they develop. Through this data analysis, students are able to identify independent and dependent variables, and make predictions from their data sets. "
$ws.Range("AT2").Value = "Improved language arts skills."
$ws.Range("AT3").Value = "Students are asked to journal and reflect on each day in the field The journal reflections are collected by the teachers and by the K12 Program Coordiantor for grading and assessment. "
$ws.Range("AU2").Value = "Better comprehension of social studies content."
$ws.Range("AU3").Value = "MOSS curriuclum focuses on holistic understanding of socio-ecological systems, including the social, economic, cultural, and scienctific connections of local and regional issues (eg: water resources in a changing climate)"
$ws.Range("AV2").Value = "Accessibility to students of all abilities and learning styles"
$ws.Range("AV3").Value = "We work to accommodate individual student needs, though we rely on the school to provide specifcs regarding IEPs or instructional aids. "

# --- Remaining "X" marker cells in row 3 (Z3:AE3) ---
$ws.Range("Z3").Value = "X"
$ws.Range("AA3").Value = "X"
$ws.Range("AB3").Value = "X"
$ws.Range("AC3").Value = "X"
$ws.Range("AD3").Value = "X"
$ws.Range("AE3").Value = "X"

# --- Numeric rating cells in row 3 (C3:S3) ---
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 2
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 5

# --- Apply the red-font header style (same as existing C1 header) to the new headers ---
$ws.Range("H1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("U1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("V1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("W1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("X1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("Y1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("AG1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("AH1").Font.Color = $ws.Range("C1").Font.Color
$ws.Range("AI1").Font.Color = $ws.Range("C1").Font.Color
